# Insert a new price record at row 574 of the "Poroto verde" sheet, shifting
# the existing rows 574:683 down to 575:684, then populate the new row with
# the new record's data (the rest of the table is left untouched because
# Excel's row-insert preserves all existing cell data of the following rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 574 and below down by one row.
$ws.Rows.Item(574).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A574").Value2 = 3
$ws.Range("B574").Value2 = "Femacal de La Calera"
$ws.Range("C574").Value2 = "Coquimbo"
$ws.Range("D574").Value2 = 45209
$ws.Range("D574").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E574").Value2 = 5
$ws.Range("F574").Value2 = 100112031
$ws.Range("G574").Value2 = "Poroto verde"
$ws.Range("H574").Value2 = "Sin especificar"
$ws.Range("I574").Value2 = "Primera"
$ws.Range("J574").Value2 = 45
$ws.Range("K574").Value2 = 40000
$ws.Range("L574").Value2 = 40000
$ws.Range("M574").Value2 = 40000
$ws.Range("N574").Value2 = "$/malla 25 kilos"
$ws.Range("O574").Value2 = "Provincia de Limarí"
$ws.Range("P574").Value2 = 1600
$ws.Range("Q574").Value2 = 25
$ws.Range("R574").Value2 = "Hortaliza"
